$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4660
$ws.Range("I64").Value = 2650
$ws.Range("K64").Value = 2650
$ws.Range("M64").Value = -2402

$ws.Range("H67").Value = 4660
$ws.Range("I67").Value = 2650
$ws.Range("K67").Value = 2650
$ws.Range("M67").Value = -1792

$ws.Range("H69").Value = 1543.8636
$ws.Range("J69").Value = 1507.8572
$ws.Range("L69").Value = 4523.571599999999
$ws.Range("N69").Value = -6271.571599999999

$ws.Range("H72").Value = 1543.8636
$ws.Range("J72").Value = 1507.8572
$ws.Range("L72").Value = 13570.7148
$ws.Range("N72").Value = -22306.7148

$ws.Range("H106").Value = 3164.7334
$ws.Range("I106").Value = 2272.8
$ws.Range("K106").Value = 2272.8
$ws.Range("M106").Value = -1641.8

$ws.Range("H125").Value = 587.05884
$ws.Range("I125").Value = 686.4
$ws.Range("K125").Value = 6177.599999999999
$ws.Range("M125").Value = -3717.599999999999

$ws.Range("H132").Value = 2032.6938
$ws.Range("I132").Value = 2033.375
$ws.Range("K132").Value = 6100.125
$ws.Range("M132").Value = -3570.125

$ws.Range("H133").Value = 51735.555
$ws.Range("J133").Value = 51735.555
$ws.Range("L133").Value = 51735.555
$ws.Range("N133").Value = -61855.555

$ws.Range("H137").Value = 115688.664
$ws.Range("I137").Value = 159211.94
$ws.Range("J137").Value = 4704.3
$ws.Range("K137").Value = 477635.82
$ws.Range("L137").Value = 14112.9
$ws.Range("M137").Value = -475085.82
$ws.Range("N137").Value = -19212.9

$ws.Range("H138").Value = 4073.817
$ws.Range("I138").Value = 8800
$ws.Range("J138").Value = 3865.3088
$ws.Range("K138").Value = 26400
$ws.Range("L138").Value = 11595.9264
$ws.Range("M138").Value = -21260
$ws.Range("N138").Value = -21875.9264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15754.077
$ws.Range("I32").Value = 10439.76
$ws.Range("J32").Value = 33468.465
$ws.Range("K32").Value = 10439.76
$ws.Range("L32").Value = 33468.465
$ws.Range("M32").Value = -10152.76
$ws.Range("N32").Value = -34042.465

$ws.Range("H45").Value = 3800.1
$ws.Range("I45").Value = 3679.762
$ws.Range("J45").Value = 4080.889
$ws.Range("K45").Value = 3679.762
$ws.Range("L45").Value = 4080.889
$ws.Range("M45").Value = -3302.762
$ws.Range("N45").Value = -4834.889

$ws.Range("H132").Value = 11331.107
$ws.Range("I132").Value = 2022.0869
$ws.Range("J132").Value = 54152.6
$ws.Range("K132").Value = 6066.2607
$ws.Range("L132").Value = 162457.8
$ws.Range("M132").Value = -3536.2607
$ws.Range("N132").Value = -167517.8

$ws.Range("H133").Value = 55997.25
$ws.Range("J133").Value = 69996.336
$ws.Range("L133").Value = 69996.336
$ws.Range("N133").Value = -75056.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3032.25
$ws.Range("I107").Value = 2266.5293
$ws.Range("K107").Value = 2266.5293
$ws.Range("M107").Value = -346.5293000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.666668
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 50
$ws.Range("L7").Value = 60
$ws.Range("M7").Value = 63
$ws.Range("N7").Value = -286

$ws.Range("H132").Value = 4047.1738
$ws.Range("I132").Value = 3056.4736
$ws.Range("K132").Value = 9169.4208
$ws.Range("M132").Value = -6639.4208

$ws.Range("H141").Value = 27270
$ws.Range("J141").Value = 27270
$ws.Range("L141").Value = 27270
$ws.Range("N141").Value = -37630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8418.583000000001
$ws.Range("I2").Value = 20025.2
$ws.Range("J2").Value = 128.14285
$ws.Range("K2").Value = 120151.2
$ws.Range("L2").Value = 768.8571000000001
$ws.Range("M2").Value = -120038.2
$ws.Range("N2").Value = -994.8571000000001

$ws.Range("H20").Value = 2236.6667
$ws.Range("I20").Value = 1100
$ws.Range("J20").Value = 2464
$ws.Range("K20").Value = 3300
$ws.Range("L20").Value = 7392
$ws.Range("M20").Value = -3073
$ws.Range("N20").Value = -7846

$ws.Range("H88").Value = 9016
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9016
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 27048
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -27904

$ws.Range("H91").Value = 9016
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9016
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 27048
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -30012

$ws.Range("H114").Value = 1196
$ws.Range("I114").Value = 540
$ws.Range("J114").Value = 1633.3334
$ws.Range("K114").Value = 1620
$ws.Range("L114").Value = 4900.0002
$ws.Range("M114").Value = 1634
$ws.Range("N114").Value = -11408.0002

$ws.Range("H129").Value = 1962.7273
$ws.Range("I129").Value = 831.8
$ws.Range("J129").Value = 2905.1667
$ws.Range("K129").Value = 2495.4
$ws.Range("L129").Value = 8715.500100000001
$ws.Range("M129").Value = 2504.6
$ws.Range("N129").Value = -18715.5001

$ws.Range("H131").Value = 723.9299999999999
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 726.1919
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2178.5757
$ws.Range("N131").Value = -12258.5757
$ws.Range("M131").Value = 3540

$ws.Range("H133").Value = 7383.091
$ws.Range("I133").Value = 1550
$ws.Range("J133").Value = 7966.4
$ws.Range("K133").Value = 4650
$ws.Range("L133").Value = 23899.2
$ws.Range("M133").Value = 410
$ws.Range("N133").Value = -34019.2

$ws.Range("H137").Value = 33343036
$ws.Range("J137").Value = 33343036
$ws.Range("L137").Value = 100029108
$ws.Range("N137").Value = -100039308

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 10000000
$ws.Range("J24").Value = 10000000
$ws.Range("L24").Value = 10000000
$ws.Range("N24").Value = -10000346

$ws.Range("H97").Value = 1749.1666
$ws.Range("J97").Value = 4255.5
$ws.Range("L97").Value = 4255.5
$ws.Range("N97").Value = -5247.5

$ws.Range("H132").Value = 75303.25999999999
$ws.Range("I132").Value = 83936.08
$ws.Range("J132").Value = 53721.2
$ws.Range("K132").Value = 251808.24
$ws.Range("L132").Value = 161163.6
$ws.Range("M132").Value = -249278.24
$ws.Range("N132").Value = -166223.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1381
$ws.Range("I46").Value = 947
$ws.Range("K46").Value = 947
$ws.Range("M46").Value = -759

$ws.Range("H122").Value = 1816149.9
$ws.Range("I122").Value = 1997064.8
$ws.Range("K122").Value = 5991194.4
$ws.Range("M122").Value = -5988744.4

$ws.Range("H132").Value = 4725
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 9266.666999999999
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 27800.001
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -32860.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080

$ws.Range("H133").Value = 44715
$ws.Range("J133").Value = 44715
$ws.Range("L133").Value = 44715
$ws.Range("N133").Value = -54835
